# Update the per-algorithm "Time" columns (Q..AD) on rows 2-20 of Sheet1.
# These hold elapsed-time measurements (one column per string-distance
# algorithm run); the benchmark was re-run and produced new timings that
# mostly collapse to a single value per row (with a few columns lagging
# slightly behind, matching the pattern of the re-measured run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2:AD2").Value = 0.02636361122131348

$ws.Range("Q3:AD3").Value = 0.03515195846557617

$ws.Range("Q4:AD4").Value = 0.02929973602294922

$ws.Range("Q5:AD5").Value = 0.03518247604370117

$ws.Range("Q6:AD6").Value = 0.06842136383056641

$ws.Range("Q7:R7").Value = 0.05885457992553711
$ws.Range("S7:AD7").Value = 0.05984210968017578

$ws.Range("Q8:AB8").Value = 0.06568002700805664
$ws.Range("AC8").Value = 0.06664466857910156
$ws.Range("AD8").Value = 0.06666874885559082

$ws.Range("Q9:S9").Value = 0.1064484119415283
$ws.Range("T9").Value = 0.1074256896972656
$ws.Range("U9").Value = 0.1064484119415283
$ws.Range("V9:AD9").Value = 0.1074256896972656

$ws.Range("Q10").Value = 0.1143500804901123
$ws.Range("R10:AD10").Value = 0.1151974201202393

$ws.Range("Q11:AD11").Value = 0.1244637966156006

$ws.Range("Q12:R12").Value = 0.09470272064208984
$ws.Range("S12:T12").Value = 0.0956730842590332
$ws.Range("U12").Value = 0.09470272064208984
$ws.Range("V12:AD12").Value = 0.0956730842590332

$ws.Range("Q13:S13").Value = 0.1498477458953857
$ws.Range("T13").Value = 0.1508252620697021
$ws.Range("U13").Value = 0.1498477458953857
$ws.Range("V13:AD13").Value = 0.1508252620697021

$ws.Range("Q14:AB14").Value = 0.1576809883117676
$ws.Range("AC14:AD14").Value = 0.158677339553833

$ws.Range("Q15:AC15").Value = 0.1094305515289307
$ws.Range("AD15").Value = 0.1104245185852051

$ws.Range("Q16:V16").Value = 0.1152205467224121
$ws.Range("W16:AD16").Value = 0.1162075996398926

$ws.Range("Q17:AD17").Value = 0.192859411239624

$ws.Range("Q18:AB18").Value = 0.201521635055542
$ws.Range("AC18:AD18").Value = 0.2020206451416016

$ws.Range("Q19:AD19").Value = 0.1389584541320801

$ws.Range("Q20:Z20").Value = 0.1437268257141113
$ws.Range("AA20:AD20").Value = 0.1445889472961426
